$wb = $excel.ActiveWorkbook

# Sheet "部门情况202401" - column O (一般贷款) for rows 3-7
# These cells are stored as text (inline strings formatted like "0.00"),
# so assign them with a leading apostrophe to force a text value instead
# of letting the numeric-looking string be auto-converted to a number.
# ClearFormats() afterwards drops the implicit "quote prefix" style that
# the apostrophe entry creates, restoring the cell to the workbook's
# default (unstyled) look - matching how the source file stores them.
$wsDept = $wb.Worksheets.Item("部门情况202401")
$wsDept.Range("O3").Value = "'135665.00"
$wsDept.Range("O4").Value = "'82140.00"
$wsDept.Range("O5").Value = "'58523.34"
$wsDept.Range("O6").Value = "'8650.00"
$wsDept.Range("O7").Value = "'1000.00"
$wsDept.Range("O3:O7").ClearFormats()

# Sheet "经办人情况202401" - column I (一般贷款) for rows 7-18 and row 27
$wsAgent = $wb.Worksheets.Item("经办人情况202401")
$wsAgent.Range("I7").Value = "'32968.00"
$wsAgent.Range("I8").Value = "'48000.00"
$wsAgent.Range("I9").Value = "'800.00"
$wsAgent.Range("I10").Value = "'53897.00"
$wsAgent.Range("I11").Value = "'18800.00"
$wsAgent.Range("I12").Value = "'7900.00"
$wsAgent.Range("I13").Value = "'15940.00"
$wsAgent.Range("I14").Value = "'39500.00"
$wsAgent.Range("I15").Value = "'9500.00"
$wsAgent.Range("I16").Value = "'39523.34"
$wsAgent.Range("I17").Value = "'9500.00"
$wsAgent.Range("I18").Value = "'8650.00"
$wsAgent.Range("I27").Value = "'1000.00"
$wsAgent.Range("I7:I18").ClearFormats()
$wsAgent.Range("I27").ClearFormats()

# Sheet "供应链放还款202401" - column I (numeric) for several rows
$wsChain = $wb.Worksheets.Item("供应链放还款202401")
$wsChain.Range("I4").Value = 0
$wsChain.Range("I6").Value = 0
$wsChain.Range("I7").Value = 109
$wsChain.Range("I9").Value = 42
$wsChain.Range("I10").Value = 11
$wsChain.Range("I11").Value = 9
$wsChain.Range("I12").Value = 1
$wsChain.Range("I13").Value = 4
$wsChain.Range("I14").Value = 4
$wsChain.Range("I16").Value = 6
$wsChain.Range("I17").Value = 186
